$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "19th "
$ws.Cells.Item(2, 2).Value = 486
$ws.Cells.Item(2, 3).Value = 644
$ws.Cells.Item(2, 4).Value = 43
$ws.Cells.Item(2, 5).Value = 23

$ws.Cells.Item(3, 1).Value = "century. "
$ws.Cells.Item(3, 2).Value = 529
$ws.Cells.Item(3, 3).Value = 644
$ws.Cells.Item(3, 4).Value = 73
$ws.Cells.Item(3, 5).Value = 23

$ws.Cells.Item(4, 1).Value = "“Romantic "
$ws.Cells.Item(4, 2).Value = 752
$ws.Cells.Item(4, 3).Value = 644
$ws.Cells.Item(4, 4).Value = 93
$ws.Cells.Item(4, 5).Value = 23

$ws.Cells.Item(5, 1).Value = "period” "
$ws.Cells.Item(5, 2).Value = 845
$ws.Cells.Item(5, 3).Value = 644
$ws.Cells.Item(5, 4).Value = 66
$ws.Cells.Item(5, 5).Value = 23

$ws.Cells.Item(6, 1).Value = "the "
$ws.Cells.Item(6, 2).Value = 656
$ws.Cells.Item(6, 3).Value = 669
$ws.Cells.Item(6, 4).Value = 32
$ws.Cells.Item(6, 5).Value = 23

$ws.Cells.Item(7, 1).Value = "Classical "
$ws.Cells.Item(7, 2).Value = 688
$ws.Cells.Item(7, 3).Value = 669
$ws.Cells.Item(7, 4).Value = 83
$ws.Cells.Item(7, 5).Value = 23

$ws.Cells.Item(8, 1).Value = "period, "
$ws.Cells.Item(8, 2).Value = 771
$ws.Cells.Item(8, 3).Value = 669
$ws.Cells.Item(8, 4).Value = 64
$ws.Cells.Item(8, 5).Value = 23

$ws.Cells.Item(9, 1).Value = "18th-century "
$ws.Cells.Item(9, 2).Value = 292
$ws.Cells.Item(9, 3).Value = 735.4
$ws.Cells.Item(9, 4).Value = 112
$ws.Cells.Item(9, 5).Value = 23

$ws.Cells.Item(10, 1).Value = "Pastoral,” "
$ws.Cells.Item(10, 2).Value = 455
$ws.Cells.Item(10, 3).Value = 1034.6
$ws.Cells.Item(10, 4).Value = 88
$ws.Cells.Item(10, 5).Value = 23

$ws.Cells.Item(11, 1).Value = "sea "
$ws.Cells.Item(11, 2).Value = 272
$ws.Cells.Item(11, 3).Value = 1126
$ws.Cells.Item(11, 4).Value = 37
$ws.Cells.Item(11, 5).Value = 23

$ws.Cells.Item(12, 1).Value = "coming "
$ws.Cells.Item(12, 2).Value = 309
$ws.Cells.Item(12, 3).Value = 1126
$ws.Cells.Item(12, 4).Value = 68
$ws.Cells.Item(12, 5).Value = 23

$ws.Cells.Item(13, 1).Value = "into "
$ws.Cells.Item(13, 2).Value = 377
$ws.Cells.Item(13, 3).Value = 1126
$ws.Cells.Item(13, 4).Value = 37
$ws.Cells.Item(13, 5).Value = 23

$ws.Cells.Item(14, 1).Value = "Fingal’s "
$ws.Cells.Item(14, 2).Value = 414
$ws.Cells.Item(14, 3).Value = 1126
$ws.Cells.Item(14, 4).Value = 72
$ws.Cells.Item(14, 5).Value = 23

$ws.Cells.Item(15, 1).Value = "Cave "
$ws.Cells.Item(15, 2).Value = 486
$ws.Cells.Item(15, 3).Value = 1126
$ws.Cells.Item(15, 4).Value = 51
$ws.Cells.Item(15, 5).Value = 23

